$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CATEGORIA D – SIN VENTAS")

# Force text format for numeric-looking code columns so leading content/codes are preserved as text
$ws.Range("A2:A7,A9:A49").NumberFormat = "@"
$ws.Range("E2:E7,E9:E49").NumberFormat = "@"
$ws.Range("W2:W7,W9:W49").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value2 = "2707090006"
$ws.Range("B2").Value2 = "KRIBENSIS PELVICACHROMIS PULCHER"
$ws.Range("C2").Value2 = ""
$ws.Range("D2").Value2 = ""
$ws.Range("E2").Value2 = "2707"
$ws.Range("F2").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G2").Value2 = 15
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 6
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 6
$ws.Range("P2").Value2 = 92
$ws.Range("Q2").Value2 = 27
$ws.Range("R2").Value2 = 180
$ws.Range("S2").Value2 = 30
$ws.Range("T2").Value2 = "Crítico"
$ws.Range("U2").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V2").Value2 = "Compra 04/05/2025"
$ws.Range("W2").Value2 = "14"

# Row 3
$ws.Range("A3").Value2 = "2707100014"
$ws.Range("B3").Value2 = "OSCAR SURTIDO"
$ws.Range("C3").Value2 = ""
$ws.Range("D3").Value2 = ""
$ws.Range("E3").Value2 = "2707"
$ws.Range("F3").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G3").Value2 = 15
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 9
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 9
$ws.Range("P3").Value2 = 92
$ws.Range("Q3").Value2 = 92
$ws.Range("R3").Value2 = 613.33
$ws.Range("S3").Value2 = 30
$ws.Range("T3").Value2 = "Crítico"
$ws.Range("U3").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 17.89€. Prioridad máxima."
$ws.Range("V3").Value2 = "Stock inicial"
$ws.Range("W3").Value2 = "14"

# Row 4
$ws.Range("A4").Value2 = "2707130074"
$ws.Range("B4").Value2 = "PEZ HACHA MARMOL"
$ws.Range("C4").Value2 = ""
$ws.Range("D4").Value2 = ""
$ws.Range("E4").Value2 = "2707"
$ws.Range("F4").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G4").Value2 = 15
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 8
$ws.Range("M4").Value2 = 0
$ws.Range("N4").Value2 = 0
$ws.Range("O4").Value2 = 8
$ws.Range("P4").Value2 = 92
$ws.Range("Q4").Value2 = 86
$ws.Range("R4").Value2 = 573.33
$ws.Range("S4").Value2 = 30
$ws.Range("T4").Value2 = "Crítico"
$ws.Range("U4").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V4").Value2 = "Compra 06/03/2025"
$ws.Range("W4").Value2 = "14"

# Row 5
$ws.Range("A5").Value2 = "2805040002"
$ws.Range("B5").Value2 = "PLANTA ESTANQUE FLOTANTE"
$ws.Range("C5").Value2 = "UNICO"
$ws.Range("D5").Value2 = "UNICO"
$ws.Range("E5").Value2 = "2805"
$ws.Range("F5").Value2 = "PLANTAS JARDIN ACUATICO"
$ws.Range("G5").Value2 = 30
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 6
$ws.Range("M5").Value2 = 0
$ws.Range("N5").Value2 = 0
$ws.Range("O5").Value2 = 6
$ws.Range("P5").Value2 = 92
$ws.Range("Q5").Value2 = 39
$ws.Range("R5").Value2 = 130
$ws.Range("S5").Value2 = 20
$ws.Range("T5").Value2 = "Crítico"
$ws.Range("U5").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V5").Value2 = "Compra 22/04/2025"
$ws.Range("W5").Value2 = "14"

# Row 6
$ws.Range("A6").Value2 = "2104090002"
$ws.Range("B6").Value2 = "DIAMANTE BICHENOV (POEPHILA BICHENOVII)"
$ws.Range("C6").Value2 = ""
$ws.Range("D6").Value2 = ""
$ws.Range("E6").Value2 = "2104"
$ws.Range("F6").Value2 = "ANIMAL VIVO PAJARO"
$ws.Range("G6").Value2 = 30
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 0
$ws.Range("L6").Value2 = 2
$ws.Range("M6").Value2 = 0
$ws.Range("N6").Value2 = 0
$ws.Range("O6").Value2 = 2
$ws.Range("P6").Value2 = 92
$ws.Range("Q6").Value2 = 9
$ws.Range("R6").Value2 = 30
$ws.Range("S6").Value2 = 0
$ws.Range("T6").Value2 = "Crítico"
$ws.Range("U6").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V6").Value2 = "Compra 22/05/2025"
$ws.Range("W6").Value2 = "14"

# Row 7
$ws.Range("A7").Value2 = "2708020004"
$ws.Range("B7").Value2 = "TELESCOPIO ROJO 6-7CM"
$ws.Range("C7").Value2 = ""
$ws.Range("D7").Value2 = ""
$ws.Range("E7").Value2 = "2708"
$ws.Range("F7").Value2 = "PECES AGUA FRIA ACUARIOFILIA"
$ws.Range("G7").Value2 = 15
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 10
$ws.Range("M7").Value2 = 0
$ws.Range("N7").Value2 = 0
$ws.Range("O7").Value2 = 10
$ws.Range("P7").Value2 = 92
$ws.Range("Q7").Value2 = 37
$ws.Range("R7").Value2 = 246.67
$ws.Range("S7").Value2 = 30
$ws.Range("T7").Value2 = "Crítico"
$ws.Range("U7").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V7").Value2 = "Compra 24/04/2025"
$ws.Range("W7").Value2 = "14"

# Row 9
$ws.Range("A9").Value2 = "2707100040"
$ws.Range("B9").Value2 = "RAMIREZI BOLIVIANO (PAPILIOCHROMIS ALTISPINOSA)"
$ws.Range("C9").Value2 = ""
$ws.Range("D9").Value2 = ""
$ws.Range("E9").Value2 = "2707"
$ws.Range("F9").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G9").Value2 = 15
$ws.Range("H9").Value2 = 0
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 0
$ws.Range("K9").Value2 = 0
$ws.Range("L9").Value2 = 6
$ws.Range("M9").Value2 = 0
$ws.Range("N9").Value2 = 0
$ws.Range("O9").Value2 = 6
$ws.Range("P9").Value2 = 92
$ws.Range("Q9").Value2 = 92
$ws.Range("R9").Value2 = 613.33
$ws.Range("S9").Value2 = 30
$ws.Range("T9").Value2 = "Crítico"
$ws.Range("U9").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 15.41€. Prioridad máxima."
$ws.Range("V9").Value2 = "Stock inicial"
$ws.Range("W9").Value2 = "14"

# Row 10
$ws.Range("A10").Value2 = "2707130030"
$ws.Range("B10").Value2 = "PEZ CUCHILLO FANTASMA"
$ws.Range("C10").Value2 = ""
$ws.Range("D10").Value2 = ""
$ws.Range("E10").Value2 = "2707"
$ws.Range("F10").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G10").Value2 = 15
$ws.Range("H10").Value2 = 0
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 0
$ws.Range("L10").Value2 = 4
$ws.Range("M10").Value2 = 0
$ws.Range("N10").Value2 = 0
$ws.Range("O10").Value2 = 4
$ws.Range("P10").Value2 = 92
$ws.Range("Q10").Value2 = 92
$ws.Range("R10").Value2 = 613.33
$ws.Range("S10").Value2 = 30
$ws.Range("T10").Value2 = "Crítico"
$ws.Range("U10").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 18.42€. Prioridad máxima."
$ws.Range("V10").Value2 = "Stock inicial"
$ws.Range("W10").Value2 = "14"

# Row 11
$ws.Range("A11").Value2 = "2707130026"
$ws.Range("B11").Value2 = "PANGASIUS SUTCHI"
$ws.Range("C11").Value2 = "5I6"
$ws.Range("D11").Value2 = "UNICO"
$ws.Range("E11").Value2 = "2707"
$ws.Range("F11").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G11").Value2 = 15
$ws.Range("H11").Value2 = 0
$ws.Range("I11").Value2 = 0
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 0
$ws.Range("L11").Value2 = 5
$ws.Range("M11").Value2 = 0
$ws.Range("N11").Value2 = 0
$ws.Range("O11").Value2 = 5
$ws.Range("P11").Value2 = 92
$ws.Range("Q11").Value2 = 9
$ws.Range("R11").Value2 = 60
$ws.Range("S11").Value2 = 0
$ws.Range("T11").Value2 = "Crítico"
$ws.Range("U11").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V11").Value2 = "Compra 22/05/2025"
$ws.Range("W11").Value2 = "14"

# Row 12
$ws.Range("A12").Value2 = "2707050028"
$ws.Range("B12").Value2 = "BETTA HELLBOY MACHO"
$ws.Range("C12").Value2 = ""
$ws.Range("D12").Value2 = ""
$ws.Range("E12").Value2 = "2707"
$ws.Range("F12").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G12").Value2 = 15
$ws.Range("H12").Value2 = 0
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = 0
$ws.Range("L12").Value2 = 3
$ws.Range("M12").Value2 = 0
$ws.Range("N12").Value2 = 0
$ws.Range("O12").Value2 = 3
$ws.Range("P12").Value2 = 92
$ws.Range("Q12").Value2 = 92
$ws.Range("R12").Value2 = 613.33
$ws.Range("S12").Value2 = 30
$ws.Range("T12").Value2 = "Crítico"
$ws.Range("U12").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 24.42€. Prioridad máxima."
$ws.Range("V12").Value2 = "Stock inicial"
$ws.Range("W12").Value2 = "14"

# Row 13
$ws.Range("A13").Value2 = "2707100024"
$ws.Range("B13").Value2 = "APISTOGRAMA BORELLI"
$ws.Range("C13").Value2 = ""
$ws.Range("D13").Value2 = ""
$ws.Range("E13").Value2 = "2707"
$ws.Range("F13").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G13").Value2 = 15
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 2
$ws.Range("M13").Value2 = 0
$ws.Range("N13").Value2 = 0
$ws.Range("O13").Value2 = 2
$ws.Range("P13").Value2 = 92
$ws.Range("Q13").Value2 = 92
$ws.Range("R13").Value2 = 613.33
$ws.Range("S13").Value2 = 30
$ws.Range("T13").Value2 = "Crítico"
$ws.Range("U13").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 3.09€. Prioridad máxima."
$ws.Range("V13").Value2 = "Stock inicial"
$ws.Range("W13").Value2 = "14"

# Row 14
$ws.Range("A14").Value2 = "2707070011"
$ws.Range("B14").Value2 = "TETRA NEON ROSA"
$ws.Range("C14").Value2 = ""
$ws.Range("D14").Value2 = ""
$ws.Range("E14").Value2 = "2707"
$ws.Range("F14").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G14").Value2 = 15
$ws.Range("H14").Value2 = 0
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 20
$ws.Range("M14").Value2 = 0
$ws.Range("N14").Value2 = 0
$ws.Range("O14").Value2 = 20
$ws.Range("P14").Value2 = 92
$ws.Range("Q14").Value2 = 24
$ws.Range("R14").Value2 = 160
$ws.Range("S14").Value2 = 30
$ws.Range("T14").Value2 = "Crítico"
$ws.Range("U14").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V14").Value2 = "Compra 07/05/2025"
$ws.Range("W14").Value2 = "14"

# Row 15
$ws.Range("A15").Value2 = "2104090011"
$ws.Range("B15").Value2 = "DIAMANTE ENMASCARADO (POEPHILA PERSONATA)"
$ws.Range("C15").Value2 = ""
$ws.Range("D15").Value2 = ""
$ws.Range("E15").Value2 = "2104"
$ws.Range("F15").Value2 = "ANIMAL VIVO PAJARO"
$ws.Range("G15").Value2 = 30
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 2
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = 0
$ws.Range("O15").Value2 = 2
$ws.Range("P15").Value2 = 92
$ws.Range("Q15").Value2 = 9
$ws.Range("R15").Value2 = 30
$ws.Range("S15").Value2 = 0
$ws.Range("T15").Value2 = "Crítico"
$ws.Range("U15").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V15").Value2 = "Compra 22/05/2025"
$ws.Range("W15").Value2 = "14"

# Row 16
$ws.Range("A16").Value2 = "2707130027"
$ws.Range("B16").Value2 = "PEZ  GLOBO AGUA DULCE"
$ws.Range("C16").Value2 = ""
$ws.Range("D16").Value2 = ""
$ws.Range("E16").Value2 = "2707"
$ws.Range("F16").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G16").Value2 = 15
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 4
$ws.Range("M16").Value2 = 0
$ws.Range("N16").Value2 = 0
$ws.Range("O16").Value2 = 4
$ws.Range("P16").Value2 = 92
$ws.Range("Q16").Value2 = 83
$ws.Range("R16").Value2 = 553.33
$ws.Range("S16").Value2 = 30
$ws.Range("T16").Value2 = "Crítico"
$ws.Range("U16").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V16").Value2 = "Compra 09/03/2025"
$ws.Range("W16").Value2 = "14"

# Row 17
$ws.Range("A17").Value2 = "2707050001"
$ws.Range("B17").Value2 = "BESUCON HELOSTOMA TEMMINCKII"
$ws.Range("C17").Value2 = ""
$ws.Range("D17").Value2 = ""
$ws.Range("E17").Value2 = "2707"
$ws.Range("F17").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 6
$ws.Range("M17").Value2 = 0
$ws.Range("N17").Value2 = 0
$ws.Range("O17").Value2 = 6
$ws.Range("P17").Value2 = 92
$ws.Range("Q17").Value2 = 92
$ws.Range("R17").Value2 = 613.33
$ws.Range("S17").Value2 = 30
$ws.Range("T17").Value2 = "Crítico"
$ws.Range("U17").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 6.22€. Prioridad máxima."
$ws.Range("V17").Value2 = "Stock inicial"
$ws.Range("W17").Value2 = "14"

# Row 18
$ws.Range("A18").Value2 = "2806010002"
$ws.Range("B18").Value2 = "COMETA SARASA XL"
$ws.Range("C18").Value2 = "10I12"
$ws.Range("D18").Value2 = "UNICO"
$ws.Range("E18").Value2 = "2806"
$ws.Range("F18").Value2 = "PECES JARDIN ACUATICO"
$ws.Range("G18").Value2 = 15
$ws.Range("H18").Value2 = 0
$ws.Range("I18").Value2 = 0
$ws.Range("J18").Value2 = 0
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = 10
$ws.Range("M18").Value2 = 0
$ws.Range("N18").Value2 = 0
$ws.Range("O18").Value2 = 10
$ws.Range("P18").Value2 = 92
$ws.Range("Q18").Value2 = 41
$ws.Range("R18").Value2 = 273.33
$ws.Range("S18").Value2 = 30
$ws.Range("T18").Value2 = "Crítico"
$ws.Range("U18").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V18").Value2 = "Compra 20/04/2025"
$ws.Range("W18").Value2 = "14"

# Row 19
$ws.Range("A19").Value2 = "2708010011"
$ws.Range("B19").Value2 = "RANCHU SURTIDO 5-6"
$ws.Range("C19").Value2 = ""
$ws.Range("D19").Value2 = ""
$ws.Range("E19").Value2 = "2708"
$ws.Range("F19").Value2 = "PECES AGUA FRIA ACUARIOFILIA"
$ws.Range("G19").Value2 = 15
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 3
$ws.Range("M19").Value2 = 0
$ws.Range("N19").Value2 = 0
$ws.Range("O19").Value2 = 3
$ws.Range("P19").Value2 = 92
$ws.Range("Q19").Value2 = 92
$ws.Range("R19").Value2 = 613.33
$ws.Range("S19").Value2 = 30
$ws.Range("T19").Value2 = "Crítico"
$ws.Range("U19").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 5.15€. Prioridad máxima."
$ws.Range("V19").Value2 = "Stock inicial"
$ws.Range("W19").Value2 = "14"

# Row 20
$ws.Range("A20").Value2 = "2707050022"
$ws.Range("B20").Value2 = "BETTA KOI"
$ws.Range("C20").Value2 = ""
$ws.Range("D20").Value2 = ""
$ws.Range("E20").Value2 = "2707"
$ws.Range("F20").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G20").Value2 = 15
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("K20").Value2 = 0
$ws.Range("L20").Value2 = 4
$ws.Range("M20").Value2 = 0
$ws.Range("N20").Value2 = 0
$ws.Range("O20").Value2 = 4
$ws.Range("P20").Value2 = 92
$ws.Range("Q20").Value2 = 92
$ws.Range("R20").Value2 = 613.33
$ws.Range("S20").Value2 = 30
$ws.Range("T20").Value2 = "Crítico"
$ws.Range("U20").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 24.08€. Prioridad máxima."
$ws.Range("V20").Value2 = "Stock inicial"
$ws.Range("W20").Value2 = "14"

# Row 21
$ws.Range("A21").Value2 = "2707130020"
$ws.Range("B21").Value2 = "LABEO ZORRO VOLADOR"
$ws.Range("C21").Value2 = ""
$ws.Range("D21").Value2 = ""
$ws.Range("E21").Value2 = "2707"
$ws.Range("F21").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G21").Value2 = 15
$ws.Range("H21").Value2 = 0
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 0
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 16
$ws.Range("M21").Value2 = 0
$ws.Range("N21").Value2 = 0
$ws.Range("O21").Value2 = 16
$ws.Range("P21").Value2 = 92
$ws.Range("Q21").Value2 = 9
$ws.Range("R21").Value2 = 60
$ws.Range("S21").Value2 = 0
$ws.Range("T21").Value2 = "Crítico"
$ws.Range("U21").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V21").Value2 = "Compra 22/05/2025"
$ws.Range("W21").Value2 = "14"

# Row 22
$ws.Range("A22").Value2 = "2104050002"
$ws.Range("B22").Value2 = "CACATUA ROSEICAPILLUS (EOLOPLUS ROSEICAPILLUS)"
$ws.Range("C22").Value2 = ""
$ws.Range("D22").Value2 = ""
$ws.Range("E22").Value2 = "2104"
$ws.Range("F22").Value2 = "ANIMAL VIVO PAJARO"
$ws.Range("G22").Value2 = 30
$ws.Range("H22").Value2 = 0
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 0
$ws.Range("L22").Value2 = 1
$ws.Range("M22").Value2 = 0
$ws.Range("N22").Value2 = 0
$ws.Range("O22").Value2 = 1
$ws.Range("P22").Value2 = 92
$ws.Range("Q22").Value2 = 29
$ws.Range("R22").Value2 = 96.67
$ws.Range("S22").Value2 = 10
$ws.Range("T22").Value2 = "Crítico"
$ws.Range("U22").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 10% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V22").Value2 = "Compra 02/05/2025"
$ws.Range("W22").Value2 = "14"

# Row 23
$ws.Range("A23").Value2 = "2707100002"
$ws.Range("B23").Value2 = "ANDINOACARA PULCHER AZUL ELECTRICO 3-4"
$ws.Range("C23").Value2 = ""
$ws.Range("D23").Value2 = ""
$ws.Range("E23").Value2 = "2707"
$ws.Range("F23").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G23").Value2 = 15
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 0
$ws.Range("L23").Value2 = 6
$ws.Range("M23").Value2 = 0
$ws.Range("N23").Value2 = 0
$ws.Range("O23").Value2 = 6
$ws.Range("P23").Value2 = 92
$ws.Range("Q23").Value2 = 41
$ws.Range("R23").Value2 = 273.33
$ws.Range("S23").Value2 = 30
$ws.Range("T23").Value2 = "Crítico"
$ws.Range("U23").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V23").Value2 = "Compra 20/04/2025"
$ws.Range("W23").Value2 = "14"

# Row 24
$ws.Range("A24").Value2 = "2707090011"
$ws.Range("B24").Value2 = "CICLIDO MULTICOLOR"
$ws.Range("C24").Value2 = ""
$ws.Range("D24").Value2 = ""
$ws.Range("E24").Value2 = "2707"
$ws.Range("F24").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G24").Value2 = 15
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("K24").Value2 = 0
$ws.Range("L24").Value2 = 4
$ws.Range("M24").Value2 = 0
$ws.Range("N24").Value2 = 0
$ws.Range("O24").Value2 = 4
$ws.Range("P24").Value2 = 92
$ws.Range("Q24").Value2 = 27
$ws.Range("R24").Value2 = 180
$ws.Range("S24").Value2 = 30
$ws.Range("T24").Value2 = "Crítico"
$ws.Range("U24").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V24").Value2 = "Compra 04/05/2025"
$ws.Range("W24").Value2 = "14"

# Row 25
$ws.Range("A25").Value2 = "2204010008"
$ws.Range("B25").Value2 = "CONEJO CABEZA LEON (NANUS ORYCTOLAGUS CUNICULUS)"
$ws.Range("C25").Value2 = ""
$ws.Range("D25").Value2 = ""
$ws.Range("E25").Value2 = "2204"
$ws.Range("F25").Value2 = "ANIMAL VIVO PEQUEÑOS MAMIFEROS"
$ws.Range("G25").Value2 = 30
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("L25").Value2 = 1
$ws.Range("M25").Value2 = 0
$ws.Range("N25").Value2 = 0
$ws.Range("O25").Value2 = 1
$ws.Range("P25").Value2 = 92
$ws.Range("Q25").Value2 = 83
$ws.Range("R25").Value2 = 276.67
$ws.Range("S25").Value2 = 30
$ws.Range("T25").Value2 = "Crítico"
$ws.Range("U25").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V25").Value2 = "Compra 09/03/2025"
$ws.Range("W25").Value2 = "14"

# Row 26
$ws.Range("A26").Value2 = "2104090001"
$ws.Range("B26").Value2 = "DIAMANTE BABERO (POEPHILA ACUTICAUDA)"
$ws.Range("C26").Value2 = ""
$ws.Range("D26").Value2 = ""
$ws.Range("E26").Value2 = "2104"
$ws.Range("F26").Value2 = "ANIMAL VIVO PAJARO"
$ws.Range("G26").Value2 = 30
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 0
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = 2
$ws.Range("M26").Value2 = 0
$ws.Range("N26").Value2 = 0
$ws.Range("O26").Value2 = 2
$ws.Range("P26").Value2 = 92
$ws.Range("Q26").Value2 = 9
$ws.Range("R26").Value2 = 30
$ws.Range("S26").Value2 = 0
$ws.Range("T26").Value2 = "Crítico"
$ws.Range("U26").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V26").Value2 = "Compra 22/05/2025"
$ws.Range("W26").Value2 = "14"

# Row 27
$ws.Range("A27").Value2 = "2707190011"
$ws.Range("B27").Value2 = "PLECOSTOMUS ORO"
$ws.Range("C27").Value2 = ""
$ws.Range("D27").Value2 = ""
$ws.Range("E27").Value2 = "2707"
$ws.Range("F27").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G27").Value2 = 15
$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 0
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 12
$ws.Range("M27").Value2 = 0
$ws.Range("N27").Value2 = 0
$ws.Range("O27").Value2 = 12
$ws.Range("P27").Value2 = 92
$ws.Range("Q27").Value2 = 92
$ws.Range("R27").Value2 = 613.33
$ws.Range("S27").Value2 = 30
$ws.Range("T27").Value2 = "Crítico"
$ws.Range("U27").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 19.99€. Prioridad máxima."
$ws.Range("V27").Value2 = "Stock inicial"
$ws.Range("W27").Value2 = "14"

# Row 28
$ws.Range("A28").Value2 = "2707130027"
$ws.Range("B28").Value2 = "PEZ GLOBO AGUA DULCE"
$ws.Range("C28").Value2 = ""
$ws.Range("D28").Value2 = ""
$ws.Range("E28").Value2 = "2707"
$ws.Range("F28").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G28").Value2 = 15
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 3
$ws.Range("M28").Value2 = 0
$ws.Range("N28").Value2 = 0
$ws.Range("O28").Value2 = 3
$ws.Range("P28").Value2 = 92
$ws.Range("Q28").Value2 = 92
$ws.Range("R28").Value2 = 613.33
$ws.Range("S28").Value2 = 30
$ws.Range("T28").Value2 = "Crítico"
$ws.Range("U28").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 6.93€. Prioridad máxima."
$ws.Range("V28").Value2 = "Stock inicial"
$ws.Range("W28").Value2 = "14"

# Row 29
$ws.Range("A29").Value2 = "2707130083"
$ws.Range("B29").Value2 = "POPONDETA PASKAI"
$ws.Range("C29").Value2 = ""
$ws.Range("D29").Value2 = ""
$ws.Range("E29").Value2 = "2707"
$ws.Range("F29").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G29").Value2 = 15
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("L29").Value2 = 6
$ws.Range("M29").Value2 = 0
$ws.Range("N29").Value2 = 0
$ws.Range("O29").Value2 = 6
$ws.Range("P29").Value2 = 92
$ws.Range("Q29").Value2 = 92
$ws.Range("R29").Value2 = 613.33
$ws.Range("S29").Value2 = 30
$ws.Range("T29").Value2 = "Crítico"
$ws.Range("U29").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 12.26€. Prioridad máxima."
$ws.Range("V29").Value2 = "Stock inicial"
$ws.Range("W29").Value2 = "14"

# Row 30
$ws.Range("A30").Value2 = "2707110001"
$ws.Range("B30").Value2 = "BOTIA CARA CABALLO"
$ws.Range("C30").Value2 = ""
$ws.Range("D30").Value2 = ""
$ws.Range("E30").Value2 = "2707"
$ws.Range("F30").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G30").Value2 = 15
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 6
$ws.Range("M30").Value2 = 0
$ws.Range("N30").Value2 = 0
$ws.Range("O30").Value2 = 6
$ws.Range("P30").Value2 = 92
$ws.Range("Q30").Value2 = 6
$ws.Range("R30").Value2 = 40
$ws.Range("S30").Value2 = 0
$ws.Range("T30").Value2 = "Crítico"
$ws.Range("U30").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V30").Value2 = "Compra 25/05/2025"
$ws.Range("W30").Value2 = "14"

# Row 31
$ws.Range("A31").Value2 = "2606050001"
$ws.Range("B31").Value2 = "INSECTO HOJA (PHYLLIUN GIGANTEUM)"
$ws.Range("C31").Value2 = ""
$ws.Range("D31").Value2 = ""
$ws.Range("E31").Value2 = "2606"
$ws.Range("F31").Value2 = "ANIMAL VIVO REPTILES"
$ws.Range("G31").Value2 = 30
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 4
$ws.Range("M31").Value2 = 0
$ws.Range("N31").Value2 = 0
$ws.Range("O31").Value2 = 4
$ws.Range("P31").Value2 = 92
$ws.Range("Q31").Value2 = 6
$ws.Range("R31").Value2 = 20
$ws.Range("S31").Value2 = 0
$ws.Range("T31").Value2 = "Crítico"
$ws.Range("U31").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V31").Value2 = "Compra 25/05/2025"
$ws.Range("W31").Value2 = "14"

# Row 32
$ws.Range("A32").Value2 = "2707050026"
$ws.Range("B32").Value2 = "BETTA MACHO LONG TAIL"
$ws.Range("C32").Value2 = ""
$ws.Range("D32").Value2 = ""
$ws.Range("E32").Value2 = "2707"
$ws.Range("F32").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G32").Value2 = 15
$ws.Range("H32").Value2 = 0
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 6
$ws.Range("M32").Value2 = 0
$ws.Range("N32").Value2 = 0
$ws.Range("O32").Value2 = 6
$ws.Range("P32").Value2 = 92
$ws.Range("Q32").Value2 = 47
$ws.Range("R32").Value2 = 313.33
$ws.Range("S32").Value2 = 30
$ws.Range("T32").Value2 = "Crítico"
$ws.Range("U32").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V32").Value2 = "Compra 14/04/2025"
$ws.Range("W32").Value2 = "14"

# Row 33
$ws.Range("A33").Value2 = "2707130012"
$ws.Range("B33").Value2 = "KILLI GARDNERI"
$ws.Range("C33").Value2 = ""
$ws.Range("D33").Value2 = ""
$ws.Range("E33").Value2 = "2707"
$ws.Range("F33").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G33").Value2 = 15
$ws.Range("H33").Value2 = 0
$ws.Range("I33").Value2 = 0
$ws.Range("J33").Value2 = 0
$ws.Range("K33").Value2 = 0
$ws.Range("L33").Value2 = 10
$ws.Range("M33").Value2 = 0
$ws.Range("N33").Value2 = 0
$ws.Range("O33").Value2 = 10
$ws.Range("P33").Value2 = 92
$ws.Range("Q33").Value2 = 6
$ws.Range("R33").Value2 = 40
$ws.Range("S33").Value2 = 0
$ws.Range("T33").Value2 = "Crítico"
$ws.Range("U33").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V33").Value2 = "Compra 25/05/2025"
$ws.Range("W33").Value2 = "14"

# Row 34
$ws.Range("A34").Value2 = "2707090002"
$ws.Range("B34").Value2 = "CICLIDO TANGANICA SURTIDO"
$ws.Range("C34").Value2 = ""
$ws.Range("D34").Value2 = ""
$ws.Range("E34").Value2 = "2707"
$ws.Range("F34").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G34").Value2 = 15
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 6
$ws.Range("M34").Value2 = 0
$ws.Range("N34").Value2 = 0
$ws.Range("O34").Value2 = 6
$ws.Range("P34").Value2 = 92
$ws.Range("Q34").Value2 = 76
$ws.Range("R34").Value2 = 506.67
$ws.Range("S34").Value2 = 30
$ws.Range("T34").Value2 = "Crítico"
$ws.Range("U34").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V34").Value2 = "Compra 16/03/2025"
$ws.Range("W34").Value2 = "14"

# Row 35
$ws.Range("A35").Value2 = "2707040007"
$ws.Range("B35").Value2 = "BARBO SAWBWA"
$ws.Range("C35").Value2 = ""
$ws.Range("D35").Value2 = ""
$ws.Range("E35").Value2 = "2707"
$ws.Range("F35").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G35").Value2 = 15
$ws.Range("H35").Value2 = 0
$ws.Range("I35").Value2 = 0
$ws.Range("J35").Value2 = 0
$ws.Range("K35").Value2 = 0
$ws.Range("L35").Value2 = 10
$ws.Range("M35").Value2 = 0
$ws.Range("N35").Value2 = 0
$ws.Range("O35").Value2 = 10
$ws.Range("P35").Value2 = 92
$ws.Range("Q35").Value2 = 92
$ws.Range("R35").Value2 = 613.33
$ws.Range("S35").Value2 = 30
$ws.Range("T35").Value2 = "Crítico"
$ws.Range("U35").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 16.87€. Prioridad máxima."
$ws.Range("V35").Value2 = "Stock inicial"
$ws.Range("W35").Value2 = "14"

# Row 36
$ws.Range("A36").Value2 = "2707130086"
$ws.Range("B36").Value2 = "LABEO ALBINO"
$ws.Range("C36").Value2 = ""
$ws.Range("D36").Value2 = ""
$ws.Range("E36").Value2 = "2707"
$ws.Range("F36").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G36").Value2 = 15
$ws.Range("H36").Value2 = 0
$ws.Range("I36").Value2 = 0
$ws.Range("J36").Value2 = 0
$ws.Range("K36").Value2 = 0
$ws.Range("L36").Value2 = 4
$ws.Range("M36").Value2 = 0
$ws.Range("N36").Value2 = 0
$ws.Range("O36").Value2 = 4
$ws.Range("P36").Value2 = 92
$ws.Range("Q36").Value2 = 76
$ws.Range("R36").Value2 = 506.67
$ws.Range("S36").Value2 = 30
$ws.Range("T36").Value2 = "Crítico"
$ws.Range("U36").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V36").Value2 = "Compra 16/03/2025"
$ws.Range("W36").Value2 = "14"

# Row 37
$ws.Range("A37").Value2 = "2708010010"
$ws.Range("B37").Value2 = "FANTAIL SURTIDO"
$ws.Range("C37").Value2 = "5CM"
$ws.Range("D37").Value2 = "UNICO"
$ws.Range("E37").Value2 = "2708"
$ws.Range("F37").Value2 = "PECES AGUA FRIA ACUARIOFILIA"
$ws.Range("G37").Value2 = 15
$ws.Range("H37").Value2 = 0
$ws.Range("I37").Value2 = 0
$ws.Range("J37").Value2 = 0
$ws.Range("K37").Value2 = 0
$ws.Range("L37").Value2 = 40
$ws.Range("M37").Value2 = 0
$ws.Range("N37").Value2 = 0
$ws.Range("O37").Value2 = 40
$ws.Range("P37").Value2 = 92
$ws.Range("Q37").Value2 = 92
$ws.Range("R37").Value2 = 613.33
$ws.Range("S37").Value2 = 30
$ws.Range("T37").Value2 = "Crítico"
$ws.Range("U37").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 35.0€. Prioridad máxima."
$ws.Range("V37").Value2 = "Stock inicial"
$ws.Range("W37").Value2 = "14"

# Row 38
$ws.Range("A38").Value2 = "2806030003"
$ws.Range("B38").Value2 = "ESTURION ALBINO (CITES CZ81A01472)"
$ws.Range("C38").Value2 = ""
$ws.Range("D38").Value2 = ""
$ws.Range("E38").Value2 = "2806"
$ws.Range("F38").Value2 = "PECES JARDIN ACUATICO"
$ws.Range("G38").Value2 = 15
$ws.Range("H38").Value2 = 0
$ws.Range("I38").Value2 = 0
$ws.Range("J38").Value2 = 0
$ws.Range("K38").Value2 = 0
$ws.Range("L38").Value2 = 8
$ws.Range("M38").Value2 = 0
$ws.Range("N38").Value2 = 0
$ws.Range("O38").Value2 = 8
$ws.Range("P38").Value2 = 92
$ws.Range("Q38").Value2 = 92
$ws.Range("R38").Value2 = 613.33
$ws.Range("S38").Value2 = 30
$ws.Range("T38").Value2 = "Crítico"
$ws.Range("U38").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 62.22€. Prioridad máxima."
$ws.Range("V38").Value2 = "Stock inicial"
$ws.Range("W38").Value2 = "14"

# Row 39
$ws.Range("A39").Value2 = "2707170006"
$ws.Range("B39").Value2 = "BADIS ESCARLATA"
$ws.Range("C39").Value2 = ""
$ws.Range("D39").Value2 = ""
$ws.Range("E39").Value2 = "2707"
$ws.Range("F39").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G39").Value2 = 15
$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
$ws.Range("J39").Value2 = 0
$ws.Range("K39").Value2 = 0
$ws.Range("L39").Value2 = 5
$ws.Range("M39").Value2 = 0
$ws.Range("N39").Value2 = 0
$ws.Range("O39").Value2 = 5
$ws.Range("P39").Value2 = 92
$ws.Range("Q39").Value2 = 27
$ws.Range("R39").Value2 = 180
$ws.Range("S39").Value2 = 30
$ws.Range("T39").Value2 = "Crítico"
$ws.Range("U39").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V39").Value2 = "Compra 04/05/2025"
$ws.Range("W39").Value2 = "14"

# Row 40
$ws.Range("A40").Value2 = "2707050018"
$ws.Range("B40").Value2 = "GURAMI SURTIDO"
$ws.Range("C40").Value2 = ""
$ws.Range("D40").Value2 = ""
$ws.Range("E40").Value2 = "2707"
$ws.Range("F40").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G40").Value2 = 15
$ws.Range("H40").Value2 = 0
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 0
$ws.Range("L40").Value2 = 15
$ws.Range("M40").Value2 = 0
$ws.Range("N40").Value2 = 0
$ws.Range("O40").Value2 = 15
$ws.Range("P40").Value2 = 92
$ws.Range("Q40").Value2 = 92
$ws.Range("R40").Value2 = 613.33
$ws.Range("S40").Value2 = 30
$ws.Range("T40").Value2 = "Crítico"
$ws.Range("U40").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 7.98€. Prioridad máxima."
$ws.Range("V40").Value2 = "Stock inicial"
$ws.Range("W40").Value2 = "14"

# Row 41
$ws.Range("A41").Value2 = "2707130084"
$ws.Range("B41").Value2 = "PANCHAX RAYADO"
$ws.Range("C41").Value2 = ""
$ws.Range("D41").Value2 = ""
$ws.Range("E41").Value2 = "2707"
$ws.Range("F41").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G41").Value2 = 15
$ws.Range("H41").Value2 = 0
$ws.Range("I41").Value2 = 0
$ws.Range("J41").Value2 = 0
$ws.Range("K41").Value2 = 0
$ws.Range("L41").Value2 = 6
$ws.Range("M41").Value2 = 0
$ws.Range("N41").Value2 = 0
$ws.Range("O41").Value2 = 6
$ws.Range("P41").Value2 = 92
$ws.Range("Q41").Value2 = 92
$ws.Range("R41").Value2 = 613.33
$ws.Range("S41").Value2 = 30
$ws.Range("T41").Value2 = "Crítico"
$ws.Range("U41").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 9.11€. Prioridad máxima."
$ws.Range("V41").Value2 = "Stock inicial"
$ws.Range("W41").Value2 = "14"

# Row 42
$ws.Range("A42").Value2 = "2707190004"
$ws.Range("B42").Value2 = "ANCISTRUS SP ADULTO"
$ws.Range("C42").Value2 = ""
$ws.Range("D42").Value2 = ""
$ws.Range("E42").Value2 = "2707"
$ws.Range("F42").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G42").Value2 = 15
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = 0
$ws.Range("L42").Value2 = 2
$ws.Range("M42").Value2 = 0
$ws.Range("N42").Value2 = 0
$ws.Range("O42").Value2 = 2
$ws.Range("P42").Value2 = 92
$ws.Range("Q42").Value2 = 76
$ws.Range("R42").Value2 = 506.67
$ws.Range("S42").Value2 = 30
$ws.Range("T42").Value2 = "Crítico"
$ws.Range("U42").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V42").Value2 = "Compra 16/03/2025"
$ws.Range("W42").Value2 = "14"

# Row 43
$ws.Range("A43").Value2 = "2707130087"
$ws.Range("B43").Value2 = "MICRORASBORA KUBOTAI"
$ws.Range("C43").Value2 = ""
$ws.Range("D43").Value2 = ""
$ws.Range("E43").Value2 = "2707"
$ws.Range("F43").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G43").Value2 = 15
$ws.Range("H43").Value2 = 0
$ws.Range("I43").Value2 = 0
$ws.Range("J43").Value2 = 0
$ws.Range("K43").Value2 = 0
$ws.Range("L43").Value2 = 10
$ws.Range("M43").Value2 = 0
$ws.Range("N43").Value2 = 0
$ws.Range("O43").Value2 = 10
$ws.Range("P43").Value2 = 92
$ws.Range("Q43").Value2 = 76
$ws.Range("R43").Value2 = 506.67
$ws.Range("S43").Value2 = 30
$ws.Range("T43").Value2 = "Crítico"
$ws.Range("U43").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V43").Value2 = "Compra 16/03/2025"
$ws.Range("W43").Value2 = "14"

# Row 44
$ws.Range("A44").Value2 = "2806030002"
$ws.Range("B44").Value2 = "ESTURION (CITES CZ18737749)"
$ws.Range("C44").Value2 = ""
$ws.Range("D44").Value2 = ""
$ws.Range("E44").Value2 = "2806"
$ws.Range("F44").Value2 = "PECES JARDIN ACUATICO"
$ws.Range("G44").Value2 = 15
$ws.Range("H44").Value2 = 0
$ws.Range("I44").Value2 = 0
$ws.Range("J44").Value2 = 0
$ws.Range("K44").Value2 = 0
$ws.Range("L44").Value2 = 4
$ws.Range("M44").Value2 = 0
$ws.Range("N44").Value2 = 0
$ws.Range("O44").Value2 = 4
$ws.Range("P44").Value2 = 92
$ws.Range("Q44").Value2 = 92
$ws.Range("R44").Value2 = 613.33
$ws.Range("S44").Value2 = 30
$ws.Range("T44").Value2 = "Crítico"
$ws.Range("U44").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 25.45€. Prioridad máxima."
$ws.Range("V44").Value2 = "Stock inicial"
$ws.Range("W44").Value2 = "14"

# Row 45
$ws.Range("A45").Value2 = "2707130055"
$ws.Range("B45").Value2 = "PEZ MARIPOSA PANTODON"
$ws.Range("C45").Value2 = ""
$ws.Range("D45").Value2 = ""
$ws.Range("E45").Value2 = "2707"
$ws.Range("F45").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G45").Value2 = 15
$ws.Range("H45").Value2 = 0
$ws.Range("I45").Value2 = 0
$ws.Range("J45").Value2 = 0
$ws.Range("K45").Value2 = 0
$ws.Range("L45").Value2 = 4
$ws.Range("M45").Value2 = 0
$ws.Range("N45").Value2 = 0
$ws.Range("O45").Value2 = 4
$ws.Range("P45").Value2 = 92
$ws.Range("Q45").Value2 = 6
$ws.Range("R45").Value2 = 40
$ws.Range("S45").Value2 = 0
$ws.Range("T45").Value2 = "Crítico"
$ws.Range("U45").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V45").Value2 = "Compra 25/05/2025"
$ws.Range("W45").Value2 = "14"

# Row 46
$ws.Range("A46").Value2 = "2806030001"
$ws.Range("B46").Value2 = "COMETA SURTIDO XL"
$ws.Range("C46").Value2 = "15CM"
$ws.Range("D46").Value2 = "UNICO"
$ws.Range("E46").Value2 = "2806"
$ws.Range("F46").Value2 = "PECES JARDIN ACUATICO"
$ws.Range("G46").Value2 = 15
$ws.Range("H46").Value2 = 0
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 4
$ws.Range("M46").Value2 = 0
$ws.Range("N46").Value2 = 0
$ws.Range("O46").Value2 = 4
$ws.Range("P46").Value2 = 92
$ws.Range("Q46").Value2 = 92
$ws.Range("R46").Value2 = 613.33
$ws.Range("S46").Value2 = 30
$ws.Range("T46").Value2 = "Crítico"
$ws.Range("U46").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 8.23€. Prioridad máxima."
$ws.Range("V46").Value2 = "Stock inicial"
$ws.Range("W46").Value2 = "14"

# Row 47
$ws.Range("A47").Value2 = "2708010014"
$ws.Range("B47").Value2 = "RYUKIN SURTIDO"
$ws.Range("C47").Value2 = "5I6"
$ws.Range("D47").Value2 = "UNICO"
$ws.Range("E47").Value2 = "2708"
$ws.Range("F47").Value2 = "PECES AGUA FRIA ACUARIOFILIA"
$ws.Range("G47").Value2 = 15
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("K47").Value2 = 0
$ws.Range("L47").Value2 = 8
$ws.Range("M47").Value2 = 0
$ws.Range("N47").Value2 = 0
$ws.Range("O47").Value2 = 8
$ws.Range("P47").Value2 = 92
$ws.Range("Q47").Value2 = 92
$ws.Range("R47").Value2 = 613.33
$ws.Range("S47").Value2 = 30
$ws.Range("T47").Value2 = "Crítico"
$ws.Range("U47").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 13.44€. Prioridad máxima."
$ws.Range("V47").Value2 = "Stock inicial"
$ws.Range("W47").Value2 = "14"

# Row 48
$ws.Range("A48").Value2 = "2707120009"
$ws.Range("B48").Value2 = "LORICARIA FILAMENTOSA"
$ws.Range("C48").Value2 = ""
$ws.Range("D48").Value2 = ""
$ws.Range("E48").Value2 = "2707"
$ws.Range("F48").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G48").Value2 = 15
$ws.Range("H48").Value2 = 0
$ws.Range("I48").Value2 = 0
$ws.Range("J48").Value2 = 0
$ws.Range("K48").Value2 = 0
$ws.Range("L48").Value2 = 4
$ws.Range("M48").Value2 = 0
$ws.Range("N48").Value2 = 0
$ws.Range("O48").Value2 = 4
$ws.Range("P48").Value2 = 92
$ws.Range("Q48").Value2 = 6
$ws.Range("R48").Value2 = 40
$ws.Range("S48").Value2 = 0
$ws.Range("T48").Value2 = "Crítico"
$ws.Range("U48").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima."
$ws.Range("V48").Value2 = "Compra 25/05/2025"
$ws.Range("W48").Value2 = "14"

# Row 49
$ws.Range("A49").Value2 = "2707110003"
$ws.Range("B49").Value2 = "BOTIA LOHACHATA"
$ws.Range("C49").Value2 = ""
$ws.Range("D49").Value2 = ""
$ws.Range("E49").Value2 = "2707"
$ws.Range("F49").Value2 = "PECES AGUA CALIENTE ACUARIOFILIA"
$ws.Range("G49").Value2 = 15
$ws.Range("H49").Value2 = 0
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 0
$ws.Range("K49").Value2 = 0
$ws.Range("L49").Value2 = 7
$ws.Range("M49").Value2 = 0
$ws.Range("N49").Value2 = 0
$ws.Range("O49").Value2 = 7
$ws.Range("P49").Value2 = 92
$ws.Range("Q49").Value2 = 92
$ws.Range("R49").Value2 = 613.33
$ws.Range("S49").Value2 = 30
$ws.Range("T49").Value2 = "Crítico"
$ws.Range("U49").Value2 = "LIQUIDACIÓN URGENTE: Aplicar descuento 30% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 11.07€. Prioridad máxima."
$ws.Range("V49").Value2 = "Stock inicial"
$ws.Range("W49").Value2 = "14"
